$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the RODOLFO row (account 004213929) entirely, shifting rows up
$ws.Rows.Item(3).Delete()

# Insert a new row after the GIANLUCA row (account 004839302), now at row 15
# after the deletion above, and populate it with the new account data.
$ws.Rows.Item(16).Insert()
$ws.Cells.Item(16, 1).NumberFormat = "@"
$ws.Cells.Item(16, 1).Value = "004165558"
$ws.Cells.Item(16, 2).Value = "GABRIEL"
$ws.Cells.Item(16, 3).Value = 200
